$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "prevalence c429" row's filepath from the placeholder
# "to be determined" to the forecasted prevalence NetCDF path.
$ws.Range("E2").Value = "/ihme/costeffectiveness/vivarium_csu_cancer/429_ets_prevalence_beta_8_phi_89.nc"

# Move the active selection to E7, matching the saved view state.
$ws.Range("E7").Select()

# Widen column E to fit the new, longer filepath text.
$ws.Columns("E").ColumnWidth = 81.5
